$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: quality_comparison
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# C1/D1: give the two blank header cells a border (top+bottom for C1,
# right+top+bottom for D1) to match the new cellXfs entries.
$ws1.Range("C1").Borders.Item(8).LineStyle = 1
$ws1.Range("C1").Borders.Item(9).LineStyle = 1

$ws1.Range("D1").Borders.Item(8).LineStyle = 1
$ws1.Range("D1").Borders.Item(9).LineStyle = 1
$ws1.Range("D1").Borders.Item(10).LineStyle = 1

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# Normalize -0 values to 0
$ws1.Range("D4").Value = 0
$ws1.Range("D5").Value = 0
$ws1.Range("D12").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: computational_comparison
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C1").Borders.Item(8).LineStyle = 1
$ws2.Range("C1").Borders.Item(9).LineStyle = 1

$ws2.Range("D1").Borders.Item(8).LineStyle = 1
$ws2.Range("D1").Borders.Item(9).LineStyle = 1
$ws2.Range("D1").Borders.Item(10).LineStyle = 1

$ws2.Range("F1").Borders.Item(8).LineStyle = 1
$ws2.Range("F1").Borders.Item(9).LineStyle = 1

$ws2.Range("G1").Borders.Item(8).LineStyle = 1
$ws2.Range("G1").Borders.Item(9).LineStyle = 1
$ws2.Range("G1").Borders.Item(10).LineStyle = 1

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string placeholder in G5
$ws2.Range("G5").ClearContents()

Write-Output "done"
